$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# --- Price/Volume(1h) updates (numeric-looking text, rows 2-15) ---
Set-TextValue $ws.Range("D2") '303.58'
Set-TextValue $ws.Range("E2") '5.21%'
Set-TextValue $ws.Range("D3") '34.90'
Set-TextValue $ws.Range("E3") '12.55%'
Set-TextValue $ws.Range("D4") '5.124'
Set-TextValue $ws.Range("E4") '4.17%'
Set-TextValue $ws.Range("D5") '0.07767'
Set-TextValue $ws.Range("E5") '5.29%'
Set-TextValue $ws.Range("D6") '2.362'
Set-TextValue $ws.Range("E6") '6.72%'
Set-TextValue $ws.Range("D7") '8.017'
Set-TextValue $ws.Range("E7") '4.18%'
Set-TextValue $ws.Range("D8") '3.939'
Set-TextValue $ws.Range("E8") '5.51%'
Set-TextValue $ws.Range("D9") '0.9272'
Set-TextValue $ws.Range("E9") '2.05%'
Set-TextValue $ws.Range("D10") '0.1015'
Set-TextValue $ws.Range("E10") '16.01%'
Set-TextValue $ws.Range("D11") '0.1795'
Set-TextValue $ws.Range("E11") '6.48%'
Set-TextValue $ws.Range("D12") '0.08557'
Set-TextValue $ws.Range("E12") '4.97%'
Set-TextValue $ws.Range("D13") '0.03310'
Set-TextValue $ws.Range("E13") '6.29%'
Set-TextValue $ws.Range("D14") '0.09890'
Set-TextValue $ws.Range("E14") '-0.59%'
Set-TextValue $ws.Range("D15") '0.001496'
Set-TextValue $ws.Range("E15") '-0.28%'

# --- Coin / Link shifts + Price/Volume(1h) updates, rows 16-23 ---
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D16") '0.005760'
Set-TextValue $ws.Range("E16") '-1.05%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D17") '3.467'
Set-TextValue $ws.Range("E17") '-0.67%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D18") '2.164'
Set-TextValue $ws.Range("E18") '4.79%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range("D19") '0.3367'
Set-TextValue $ws.Range("E19") '1.17%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range("D20") '0.1309'
Set-TextValue $ws.Range("E20") '1.05%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range("D21") '4.311'
Set-TextValue $ws.Range("E21") '12.64%'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range("D22") '0.2385'
Set-TextValue $ws.Range("E22") '12.38%'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range("D23") '0.04562'
Set-TextValue $ws.Range("E23") '0.19%'

# --- Price/Volume(1h) updates, rows 24-27 ---
Set-TextValue $ws.Range("D24") '0.001217'
Set-TextValue $ws.Range("E24") '0.54%'
Set-TextValue $ws.Range("D25") '0.004460'
Set-TextValue $ws.Range("E25") '7.61%'
Set-TextValue $ws.Range("D26") '0.0001250'
Set-TextValue $ws.Range("E26") '-3.95%'
Set-TextValue $ws.Range("D27") '0.0003697'
Set-TextValue $ws.Range("E27") '8.78%'

# --- Price/Volume(1h) updates, rows 39-51 ---
Set-TextValue $ws.Range("D39") '0.01789'
Set-TextValue $ws.Range("E39") '13.16%'
Set-TextValue $ws.Range("D40") '0.04765'
Set-TextValue $ws.Range("E40") '6.66%'
Set-TextValue $ws.Range("D41") '0.007734'
Set-TextValue $ws.Range("E41") '5.25%'
Set-TextValue $ws.Range("D42") '0.1412'
Set-TextValue $ws.Range("E42") '6.62%'
Set-TextValue $ws.Range("D43") '0.007089'
Set-TextValue $ws.Range("E43") '-25.80%'
Set-TextValue $ws.Range("D44") '0.002106'
Set-TextValue $ws.Range("E44") '-5.25%'
Set-TextValue $ws.Range("D45") '0.009538'
Set-TextValue $ws.Range("E45") '13.20%'
Set-TextValue $ws.Range("D46") '0.00006115'
Set-TextValue $ws.Range("E46") '0.08%'
Set-TextValue $ws.Range("D47") '0.00000000749'
Set-TextValue $ws.Range("E47") '-0.19%'
Set-TextValue $ws.Range("D48") '2.736'
Set-TextValue $ws.Range("E48") '29.97%'
Set-TextValue $ws.Range("D49") '0.001998'
Set-TextValue $ws.Range("E49") '-0.20%'
Set-TextValue $ws.Range("D50") '0.00002098'
Set-TextValue $ws.Range("E50") '-0.19%'
Set-TextValue $ws.Range("D51") '0.0001998'
Set-TextValue $ws.Range("E51") '-0.19%'
